$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the missing ORDER BY / LIMIT clause to the Neo4j query in B2 ---
# (Doing this as an append to the existing value avoids re-typing the whole
#  Cypher query, which is full of characters that are special to PowerShell.)
$suffix = "`n order By ss.study_subject_id ASC LIMIT 100 "
$currentQuery = $ws.Range("B2").Value()
$ws.Range("B2").Value = $currentQuery + $suffix

# --- Row 2 now wraps one extra line, so grow its height to match ---
$ws.Rows.Item(2).RowHeight = 331.2

# --- Scroll the sheet view back to the top row (topLeftCell B2 -> B1) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("C2").Select()
